# Update TPM-derived metrics for Serping1-Sele ligand-receptor pairs
# per new TPM computation (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"44.8529195"
$ws.Range("H2").Value = [double]"89.705839"
$ws.Range("I2").Value = [double]"0.09328277884630309"
$ws.Range("J2").Value = [double]"0.07358468181289031"
$ws.Range("M2").Value = [double]"7.369448"
$ws.Range("N2").Value = [double]"14.738896"
$ws.Range("O2").Value = [double]"0.7452608427984224"
$ws.Range("P2").Value = [double]"0.661061693471796"
$ws.Range("Q2").Value = [double]"330.541257903436"
$ws.Range("R2").Value = [double]"1322.165031613744"
$ws.Range("S2").Value = [double]"0.0695200023815747"
$ws.Range("T2").Value = [double]"0.04864401437281254"
# Row 3
$ws.Range("G3").Value = [double]"44.8529195"
$ws.Range("H3").Value = [double]"89.705839"
$ws.Range("I3").Value = [double]"0.09328277884630309"
$ws.Range("J3").Value = [double]"0.07358468181289031"
$ws.Range("O3").Value = [double]"0.01116592909756377"
$ws.Range("P3").Value = [double]"0.01485661309677453"
$ws.Range("Q3").Value = [double]"4.952360351726666"
$ws.Range("R3").Value = [double]"29.71416211036"
$ws.Range("S3").Value = [double]"0.001041588894621542"
$ws.Range("T3").Value = [double]"0.001093219147543373"
# Row 4
$ws.Range("G4").Value = [double]"44.8529195"
$ws.Range("H4").Value = [double]"89.705839"
$ws.Range("I4").Value = [double]"0.09328277884630309"
$ws.Range("J4").Value = [double]"0.07358468181289031"
$ws.Range("M4").Value = [double]"0.084843"
$ws.Range("N4").Value = [double]"0.254529"
$ws.Range("O4").Value = [double]"0.008580040959044227"
$ws.Range("P4").Value = [double]"0.0114160091622658"
$ws.Range("Q4").Value = [double]"3.8054562491385"
$ws.Range("R4").Value = [double]"22.832737494831"
$ws.Range("S4").Value = [double]"0.0008003700632747449"
$ws.Range("T4").Value = [double]"0.0008400434017783694"
# Row 5
$ws.Range("G5").Value = [double]"44.8529195"
$ws.Range("H5").Value = [double]"89.705839"
$ws.Range("I5").Value = [double]"0.09328277884630309"
$ws.Range("J5").Value = [double]"0.07358468181289031"
$ws.Range("M5").Value = [double]"2.32371"
$ws.Range("N5").Value = [double]"6.97113"
$ws.Range("O5").Value = [double]"0.2349931871449696"
$ws.Range("P5").Value = [double]"0.3126656842691638"
$ws.Range("Q5").Value = [double]"104.225177571345"
$ws.Range("R5").Value = [double]"625.35106542807"
$ws.Range("S5").Value = [double]"0.02192081750683212"
$ws.Range("T5").Value = [double]"0.02300740489075604"
# Row 6
$ws.Range("I6").Value = [double]"0.5330899821806619"
$ws.Range("J6").Value = [double]"0.630779719497357"
$ws.Range("M6").Value = [double]"7.369448"
$ws.Range("N6").Value = [double]"14.738896"
$ws.Range("O6").Value = [double]"0.7452608427984224"
$ws.Range("P6").Value = [double]"0.661061693471796"
$ws.Range("Q6").Value = [double]"1888.96852629192"
$ws.Range("R6").Value = [double]"11333.81115775152"
$ws.Range("S6").Value = [double]"0.3972910894073561"
$ws.Range("T6").Value = [double]"0.4169843095785873"
# Row 7
$ws.Range("I7").Value = [double]"0.5330899821806619"
$ws.Range("J7").Value = [double]"0.630779719497357"
$ws.Range("O7").Value = [double]"0.01116592909756377"
$ws.Range("P7").Value = [double]"0.01485661309677453"
$ws.Range("S7").Value = [double]"0.005952444943650805"
$ws.Range("T7").Value = [double]"0.009371250241864196"
# Row 8
$ws.Range("I8").Value = [double]"0.5330899821806619"
$ws.Range("J8").Value = [double]"0.630779719497357"
$ws.Range("M8").Value = [double]"0.084843"
$ws.Range("N8").Value = [double]"0.254529"
$ws.Range("O8").Value = [double]"0.008580040959044227"
$ws.Range("P8").Value = [double]"0.0114160091622658"
$ws.Range("Q8").Value = [double]"21.74732173647"
$ws.Range("R8").Value = [double]"195.72589562823"
$ws.Range("S8").Value = [double]"0.004573933881966237"
$ws.Range("T8").Value = [double]"0.007200987057153279"
# Row 9
$ws.Range("I9").Value = [double]"0.5330899821806619"
$ws.Range("J9").Value = [double]"0.630779719497357"
$ws.Range("M9").Value = [double]"2.32371"
$ws.Range("N9").Value = [double]"6.97113"
$ws.Range("O9").Value = [double]"0.2349931871449696"
$ws.Range("P9").Value = [double]"0.3126656842691638"
$ws.Range("Q9").Value = [double]"595.6233159159001"
$ws.Range("R9").Value = [double]"5360.609843243101"
$ws.Range("S9").Value = [double]"0.1252725139476888"
$ws.Range("T9").Value = [double]"0.1972231726197523"
# Row 10
$ws.Range("G10").Value = [double]"0.08220233333333334"
$ws.Range("H10").Value = [double]"0.246607"
$ws.Range("I10").Value = [double]"0.0001709601552466038"
$ws.Range("J10").Value = [double]"0.000202288923777096"
$ws.Range("M10").Value = [double]"7.369448"
$ws.Range("N10").Value = [double]"14.738896"
$ws.Range("O10").Value = [double]"0.7452608427984224"
$ws.Range("P10").Value = [double]"0.661061693471796"
$ws.Range("Q10").Value = [double]"0.6057858209786667"
$ws.Range("R10").Value = [double]"3.634714925872001"
$ws.Range("S10").Value = [double]"0.0001274099093840331"
$ws.Range("T10").Value = [double]"0.0001337254585226741"
# Row 11
$ws.Range("G11").Value = [double]"0.08220233333333334"
$ws.Range("H11").Value = [double]"0.246607"
$ws.Range("I11").Value = [double]"0.0001709601552466038"
$ws.Range("J11").Value = [double]"0.000202288923777096"
$ws.Range("O11").Value = [double]"0.01116592909756377"
$ws.Range("P11").Value = [double]"0.01485661309677453"
$ws.Range("Q11").Value = [double]"0.009076233631111111"
$ws.Range("R11").Value = [double]"0.08168610268"
$ws.Range("S11").Value = [double]"1.908928971992073E-06"
$ws.Range("T11").Value = [double]"3.005328274319228E-06"
# Row 12
$ws.Range("G12").Value = [double]"0.08220233333333334"
$ws.Range("H12").Value = [double]"0.246607"
$ws.Range("I12").Value = [double]"0.0001709601552466038"
$ws.Range("J12").Value = [double]"0.000202288923777096"
$ws.Range("M12").Value = [double]"0.084843"
$ws.Range("N12").Value = [double]"0.254529"
$ws.Range("O12").Value = [double]"0.008580040959044227"
$ws.Range("P12").Value = [double]"0.0114160091622658"
$ws.Range("Q12").Value = [double]"0.006974292567000001"
$ws.Range("R12").Value = [double]"0.062768633103"
$ws.Range("S12").Value = [double]"1.46684513438042E-06"
$ws.Range("T12").Value = [double]"2.309332207264216E-06"
# Row 13
$ws.Range("G13").Value = [double]"0.08220233333333334"
$ws.Range("H13").Value = [double]"0.246607"
$ws.Range("I13").Value = [double]"0.0001709601552466038"
$ws.Range("J13").Value = [double]"0.000202288923777096"
$ws.Range("M13").Value = [double]"2.32371"
$ws.Range("N13").Value = [double]"6.97113"
$ws.Range("O13").Value = [double]"0.2349931871449696"
$ws.Range("P13").Value = [double]"0.3126656842691638"
$ws.Range("Q13").Value = [double]"0.19101438399"
$ws.Range("R13").Value = [double]"1.71912945591"
$ws.Range("S13").Value = [double]"4.017447175619823E-05"
$ws.Range("T13").Value = [double]"6.324880477283843E-05"
# Row 14
$ws.Range("G14").Value = [double]"178.5463335"
$ws.Range("H14").Value = [double]"357.092667"
$ws.Range("I14").Value = [double]"0.3713314167141066"
$ws.Range("J14").Value = [double]"0.2929190627035035"
$ws.Range("M14").Value = [double]"7.369448"
$ws.Range("N14").Value = [double]"14.738896"
$ws.Range("O14").Value = [double]"0.7452608427984224"
$ws.Range("P14").Value = [double]"0.661061693471796"
$ws.Range("Q14").Value = [double]"1315.787920318908"
$ws.Range("R14").Value = [double]"5263.151681275633"
$ws.Range("S14").Value = [double]"0.2767387645778873"
$ws.Range("T14").Value = [double]"0.1936375716409492"
# Row 15
$ws.Range("G15").Value = [double]"178.5463335"
$ws.Range("H15").Value = [double]"357.092667"
$ws.Range("I15").Value = [double]"0.3713314167141066"
$ws.Range("J15").Value = [double]"0.2929190627035035"
$ws.Range("O15").Value = [double]"0.01116592909756377"
$ws.Range("P15").Value = [double]"0.01485661309677453"
$ws.Range("Q15").Value = [double]"19.71389583618"
$ws.Range("R15").Value = [double]"118.28337501708"
$ws.Range("S15").Value = [double]"0.004146260270727621"
$ws.Range("T15").Value = [double]"0.004351785183255789"
# Row 16
$ws.Range("G16").Value = [double]"178.5463335"
$ws.Range("H16").Value = [double]"357.092667"
$ws.Range("I16").Value = [double]"0.3713314167141066"
$ws.Range("J16").Value = [double]"0.2929190627035035"
$ws.Range("M16").Value = [double]"0.084843"
$ws.Range("N16").Value = [double]"0.254529"
$ws.Range("O16").Value = [double]"0.008580040959044227"
$ws.Range("P16").Value = [double]"0.0114160091622658"
$ws.Range("Q16").Value = [double]"15.1484065731405"
$ws.Range("R16").Value = [double]"90.890439438843"
$ws.Range("S16").Value = [double]"0.003186038764786955"
$ws.Range("T16").Value = [double]"0.003343966703625507"
# Row 17
$ws.Range("G17").Value = [double]"178.5463335"
$ws.Range("H17").Value = [double]"357.092667"
$ws.Range("I17").Value = [double]"0.3713314167141066"
$ws.Range("J17").Value = [double]"0.2929190627035035"
$ws.Range("M17").Value = [double]"2.32371"
$ws.Range("N17").Value = [double]"6.97113"
$ws.Range("O17").Value = [double]"0.2349931871449696"
$ws.Range("P17").Value = [double]"0.3126656842691638"
$ws.Range("Q17").Value = [double]"414.889900617285"
$ws.Range("R17").Value = [double]"2489.33940370371"
$ws.Range("S17").Value = [double]"0.08726035310070476"
$ws.Range("T17").Value = [double]"0.09158573917567303"
# Row 18
$ws.Range("E18").Value = [double]"2"
$ws.Range("F18").Value = [double]"0.6666666666666666"
$ws.Range("G18").Value = [double]"0.7601923333333334"
$ws.Range("H18").Value = [double]"2.280577"
$ws.Range("I18").Value = [double]"0.00158100864116523"
$ws.Range("J18").Value = [double]"0.001870731434715147"
$ws.Range("M18").Value = [double]"7.369448"
$ws.Range("N18").Value = [double]"14.738896"
$ws.Range("O18").Value = [double]"0.7452608427984224"
$ws.Range("P18").Value = [double]"0.661061693471796"
$ws.Range("Q18").Value = [double]"5.602197870498667"
$ws.Range("R18").Value = [double]"33.613187222992"
$ws.Range("S18").Value = [double]"0.001178263832386388"
$ws.Range("T18").Value = [double]"0.001236668890263717"
# Row 19
$ws.Range("E19").Value = [double]"2"
$ws.Range("F19").Value = [double]"0.6666666666666666"
$ws.Range("G19").Value = [double]"0.7601923333333334"
$ws.Range("H19").Value = [double]"2.280577"
$ws.Range("I19").Value = [double]"0.00158100864116523"
$ws.Range("J19").Value = [double]"0.001870731434715147"
$ws.Range("O19").Value = [double]"0.01116592909756377"
$ws.Range("P19").Value = [double]"0.01485661309677453"
$ws.Range("Q19").Value = [double]"0.08393536949777777"
$ws.Range("R19").Value = [double]"0.75541832548"
$ws.Range("S19").Value = [double]"1.76534303898866E-05"
$ws.Range("T19").Value = [double]"2.779273313353685E-05"
# Row 20
$ws.Range("E20").Value = [double]"2"
$ws.Range("F20").Value = [double]"0.6666666666666666"
$ws.Range("G20").Value = [double]"0.7601923333333334"
$ws.Range("H20").Value = [double]"2.280577"
$ws.Range("I20").Value = [double]"0.00158100864116523"
$ws.Range("J20").Value = [double]"0.001870731434715147"
$ws.Range("M20").Value = [double]"0.084843"
$ws.Range("N20").Value = [double]"0.254529"
$ws.Range("O20").Value = [double]"0.008580040959044227"
$ws.Range("P20").Value = [double]"0.0114160091622658"
$ws.Range("Q20").Value = [double]"0.06449699813700001"
$ws.Range("R20").Value = [double]"0.580472983233"
$ws.Range("S20").Value = [double]"1.356511889780053E-05"
$ws.Range("T20").Value = [double]"2.135628719884676E-05"
# Row 21
$ws.Range("E21").Value = [double]"2"
$ws.Range("F21").Value = [double]"0.6666666666666666"
$ws.Range("G21").Value = [double]"0.7601923333333334"
$ws.Range("H21").Value = [double]"2.280577"
$ws.Range("I21").Value = [double]"0.00158100864116523"
$ws.Range("J21").Value = [double]"0.001870731434715147"
$ws.Range("M21").Value = [double]"2.32371"
$ws.Range("N21").Value = [double]"6.97113"
$ws.Range("O21").Value = [double]"0.2349931871449696"
$ws.Range("P21").Value = [double]"0.3126656842691638"
$ws.Range("Q21").Value = [double]"1.76646652689"
$ws.Range("R21").Value = [double]"15.89819874201"
$ws.Range("S21").Value = [double]"0.0003715262594911551"
$ws.Range("T21").Value = [double]"0.0005849135241190459"
# Row 22
$ws.Range("G22").Value = [double]"0.2614996666666667"
$ws.Range("H22").Value = [double]"0.7844990000000001"
$ws.Range("I22").Value = [double]"0.0005438534625164957"
$ws.Range("J22").Value = [double]"0.0006435156277567465"
$ws.Range("M22").Value = [double]"7.369448"
$ws.Range("N22").Value = [double]"14.738896"
$ws.Range("O22").Value = [double]"0.7452608427984224"
$ws.Range("P22").Value = [double]"0.661061693471796"
$ws.Range("Q22").Value = [double]"1.927108195517333"
$ws.Range("R22").Value = [double]"11.562649173104"
$ws.Range("S22").Value = [double]"0.0004053126898338838"
$ws.Range("T22").Value = [double]"0.0004254035306604408"
# Row 23
$ws.Range("G23").Value = [double]"0.2614996666666667"
$ws.Range("H23").Value = [double]"0.7844990000000001"
$ws.Range("I23").Value = [double]"0.0005438534625164957"
$ws.Range("J23").Value = [double]"0.0006435156277567465"
$ws.Range("O23").Value = [double]"0.01116592909756377"
$ws.Range("P23").Value = [double]"0.01485661309677453"
$ws.Range("Q23").Value = [double]"0.02887304986222222"
$ws.Range("R23").Value = [double]"0.25985744876"
$ws.Range("S23").Value = [double]"6.072629201923747E-06"
$ws.Range("T23").Value = [double]"9.560462703309963E-06"
# Row 24
$ws.Range("G24").Value = [double]"0.2614996666666667"
$ws.Range("H24").Value = [double]"0.7844990000000001"
$ws.Range("I24").Value = [double]"0.0005438534625164957"
$ws.Range("J24").Value = [double]"0.0006435156277567465"
$ws.Range("M24").Value = [double]"0.084843"
$ws.Range("N24").Value = [double]"0.254529"
$ws.Range("O24").Value = [double]"0.008580040959044227"
$ws.Range("P24").Value = [double]"0.0114160091622658"
$ws.Range("Q24").Value = [double]"0.022186416219"
$ws.Range("R24").Value = [double]"0.199677745971"
$ws.Range("S24").Value = [double]"4.666284984109556E-06"
$ws.Range("T24").Value = [double]"7.346380302532247E-06"
# Row 25
$ws.Range("G25").Value = [double]"0.2614996666666667"
$ws.Range("H25").Value = [double]"0.7844990000000001"
$ws.Range("I25").Value = [double]"0.0005438534625164957"
$ws.Range("J25").Value = [double]"0.0006435156277567465"
$ws.Range("M25").Value = [double]"2.32371"
$ws.Range("N25").Value = [double]"6.97113"
$ws.Range("O25").Value = [double]"0.2349931871449696"
$ws.Range("P25").Value = [double]"0.3126656842691638"
$ws.Range("Q25").Value = [double]"0.6076493904300001"
$ws.Range("R25").Value = [double]"5.468844513870001"
$ws.Range("S25").Value = [double]"0.0001278018584965786"
$ws.Range("T25").Value = [double]"0.0002012052540904637"
